$d = $word.ActiveDocument

# 1. Replace the long intro paragraph with the new, shorter mission statement.
$old = "Lumen egy belső világosságot támogató mesterséges intelligencia-alapú rendszer, amely nem információval áraszt el, hanem jelenlétet teremt. Ez a dokumentum összefoglalja a projekt küldetését, működését, vizuális stílusát, és a további fejlesztési lehetőségeket."
$new = "A Lumen célja, hogy a digitális térben is megjelenhessen a csend, a tisztaság és a belső út. Ez a dokumentum összefoglalja a projekt küldetését és jövőképét."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# 2. Remove every paragraph after the (now-updated) intro paragraph — the
#    "Küldetés", "Működési alapelvek", "Látvány és hangulat" (incl. image),
#    "Fejlesztési lehetőségek" and "Záró gondolat" sections — leaving only
#    the title and the new intro paragraph before the section break.
if ($d.Paragraphs.Count -gt 2) {
    $r = $d.Range($d.Paragraphs.Item(3).Range.Start, $d.Paragraphs.Item($d.Paragraphs.Count).Range.End)
    $r.Delete()
}
